# Commit: "11232040 added review fa"
# Replace the "Georgia Anne Muldrow / Jyoti" tracklist (georgiaannemuldrow2)
# with the "The Alchemist / Freddie Gibbs - Alfredo" tracklist (fgibbsalchemist1),
# trimmed from 15 to 10 tracks, on Sheet1 and Sheet3 (the two sheets backing the
# web-query defined name / table), and update Sheet2's selection to match the
# shrunk range. Sheet2 is fully formula-driven off Sheet1, so it recalculates
# automatically.

$wb = $excel.ActiveWorkbook

function Set-Tracklist($ws) {
    # No. | Title | Composer(s) | Performer | Time (fraction of day)
    $ws.Range("A2").Value = 1
    $ws.Range("B2").Value = 1985
    $ws.Range("C2").Value = "Frederick Tipton, Daniel Maman"
    $ws.Range("D2").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E2").Value = 0.10555555555555556

    $ws.Range("A3").Value = 2
    $ws.Range("B3").Value = "God Is Perfect"
    $ws.Range("C3").Value = "Frederick Tipton, Daniel Maman"
    $ws.Range("D3").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E3").Value = 0.16597222222222222

    $ws.Range("A4").Value = 3
    $ws.Range("B4").Value = "Scottie Beam"
    $ws.Range("C4").Value = "Frederick Tipton, Daniel Maman, Norman Whiteside, William Roberts III"
    $ws.Range("D4").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E4").Value = 0.16944444444444443

    $ws.Range("A5").Value = 4
    $ws.Range("B5").Value = "Look at Me"
    $ws.Range("C5").Value = "Frederick Tipton, Daniel Maman"
    $ws.Range("D5").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E5").Value = 0.10625

    $ws.Range("A6").Value = 5
    $ws.Range("B6").Value = "Frank Lucas"
    $ws.Range("C6").Value = "Frederick Tipton, Daniel Maman, Jeremie Pennick"
    $ws.Range("D6").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E6").Value = 0.19513888888888889

    $ws.Range("A7").Value = 6
    $ws.Range("B7").Value = "Something to Rap About"
    $ws.Range("C7").Value = "Frederick Tipton, Daniel Maman, Nick Walker, Tyler Okonma"
    $ws.Range("D7").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E7").Value = 0.19583333333333333

    $ws.Range("A8").Value = 7
    $ws.Range("B8").Value = "Baby `$hit"
    $ws.Range("C8").Value = "Frederick Tipton, Daniel Maman"
    $ws.Range("D8").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E8").Value = 0.15

    $ws.Range("A9").Value = 8
    $ws.Range("B9").Value = "Babies & Fools"
    $ws.Range("C9").Value = "Frederick Tipton, Daniel Maman, Demond Price"
    $ws.Range("D9").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E9").Value = 0.14305555555555557

    $ws.Range("A10").Value = 9
    $ws.Range("B10").Value = "Skinny Suge"
    $ws.Range("C10").Value = "Frederick Tipton, Daniel Maman"
    $ws.Range("D10").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E10").Value = 0.11944444444444445

    $ws.Range("A11").Value = 10
    $ws.Range("B11").Value = "All Glass"
    $ws.Range("C11").Value = "Frederick Tipton, Daniel Maman"
    $ws.Range("D11").Value = "The Alchemist / Freddie Gibbs"
    $ws.Range("E11").Value = 0.10694444444444444

    # Tracks 11-15 no longer exist - drop the old data but keep formatting
    # (this matches a plain Excel "Clear Contents" on A12:E16).
    $ws.Range("A12:E16").ClearContents()

    # Column widths (Sheet1/Sheet3): B,C,D,E get new widths in the edit.
    # ColumnWidth (char units, MDW=7) = stored_width - 5/7; stored widths
    # taken from the target OOXML, rounded to the nearest 1/7 the engine
    # can actually represent.
    $ws.Columns.Item(2).ColumnWidth = 22.0
    $ws.Columns.Item(3).ColumnWidth = 64.57142857142857
    $ws.Columns.Item(4).ColumnWidth = 27.571428571428573
    $ws.Columns.Item(5).ColumnWidth = 5.0
}

Set-Tracklist($wb.Worksheets.Item("Sheet1"))
Set-Tracklist($wb.Worksheets.Item("Sheet3"))

# Defined names: georgiaannemuldrow2 -> fgibbsalchemist1, range shrinks to the
# new 10-track extent ($E$16 -> $E$11). Update RefersTo before Name to avoid
# a transient name collision between the two identically-named entries.
$n1 = $wb.Names.Item(1)
$n2 = $wb.Names.Item(2)
$n1.RefersTo = "=Sheet1!`$A`$1:`$E`$11"
$n2.RefersTo = "=Sheet3!`$A`$1:`$E`$11"
$n1.Name = "fgibbsalchemist1"
$n2.Name = "fgibbsalchemist1"

# Sheet2's selection shrinks along with the table (K19/A3:K19 -> K14/A3:K14).
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A3:K14").Select()
